# Flow_Extraction project: simplify the "Object Creation" template sheet
# down to a single "Object Name" column with one example row, clearing the
# rest of the sample data out to a blank template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused "Object IP" and "Object Groups" columns (old B:C).
# This shifts the old column D -> B and the old column G -> E, matching
# their original widths/styles.
$ws.Range("B:C").Delete()

# Replace the sample host entry and blank out the remaining example rows,
# leaving a clean single-column template.
$ws.Range("A2").Value = "Host_10.0.0.1"
$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()

# Move the active selection to the first data cell.
$ws.Range("A2").Select() | Out-Null

# Portrait page orientation for printing.
$ws.PageSetup.Orientation = 1
